# Updates to code, figures, and manuscript
# rf_performance_metrics.xlsx - refresh derived RF performance metrics

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Disable concurrent/multi-threaded calculation for this workbook
# (corresponds to calcPr concurrentCalc="0")
$excel.MultiThreadedCalculation.Enabled = $false

# Refresh computed metric values (perc.var / rmse columns)
$ws.Range("C2").Value = 24
$ws.Range("D2").Value = 7

$ws.Range("C3").Value = 27
$ws.Range("D3").Value = 4

$ws.Range("C4").Value = 22

$ws.Range("C5").Value = 18

$ws.Range("C6").Value = 17

$ws.Range("C7").Value = 24

$ws.Range("C8").Value = 20

$ws.Range("C9").Value = 11

$ws.Range("C10").Value = 41
$ws.Range("D10").Value = 0.09

$ws.Range("C13").Value = 36

# Leave the selection where the author left off editing
$ws.Range("D14").Select()
